$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (must be written as text to match
# the source workbook, which stores these as inline/shared strings, not numbers).
$updates = @{
    'D2' = '320.66'
    'E2' = '7.87%'
    'G2' = '8'
    'D3' = '49.12'
    'E3' = '17.68%'
    'G3' = '8'
    'D4' = '5.252'
    'E4' = '4.71%'
    'G4' = '8'
    'D5' = '0.08106'
    'E5' = '7.74%'
    'G5' = '8'
    'D6' = '4.624'
    'E6' = '5.65%'
    'G6' = '8'
    'D7' = '1.660'
    'E7' = '3.89%'
    'G7' = '8'
    'D8' = '1.188'
    'E8' = '29.48%'
    'G8' = '8'
    'D9' = '0.1314'
    'G9' = '8'
    'D10' = '0.1935'
    'E10' = '6.10%'
    'G10' = '8'
    'D11' = '0.09590'
    'E11' = '7.75%'
    'G11' = '8'
    'D12' = '0.04687'
    'E12' = '15.45%'
    'G12' = '8'
    'E13' = '-0.04%'
    'G13' = '8'
    'D14' = '0.001332'
    'E14' = '3.67%'
    'G14' = '8'
    'D15' = '0.005953'
    'E15' = '1.98%'
    'G15' = '8'
    'D16' = '3.366'
    'E16' = '0.79%'
    'G16' = '8'
    'D17' = '2.433'
    'E17' = '1.32%'
    'G17' = '8'
    'D18' = '0.3395'
    'E18' = '1.95%'
    'G18' = '8'
    'D19' = '8.159'
    'E19' = '-1.39%'
    'G19' = '8'
    'D20' = '0.1412'
    'E20' = '4.39%'
    'G20' = '8'
    'G21' = '8'
    'D22' = '0.04286'
    'E22' = '4.66%'
    'G22' = '8'
    'D23' = '0.001308'
    'E23' = '3.26%'
    'G23' = '8'
    'D24' = '0.004243'
    'E24' = '8.88%'
    'G24' = '8'
    'D25' = '0.0001350'
    'E25' = '3.76%'
    'G25' = '8'
    'G26' = '8'
    'G27' = '8'
    'G28' = '8'
    'G29' = '8'
    'G30' = '8'
    'G31' = '8'
    'G32' = '8'
    'G33' = '8'
    'G34' = '8'
    'G35' = '8'
    'G36' = '8'
    'G37' = '8'
    'D38' = '0.02676'
    'E38' = '11.55%'
    'G38' = '8'
    'D39' = '0.05578'
    'E39' = '7.58%'
    'G39' = '8'
    'D40' = '0.006302'
    'E40' = '-0.08%'
    'G40' = '8'
    'D41' = '0.007684'
    'E41' = '-1.30%'
    'G41' = '8'
    'D42' = '0.1439'
    'E42' = '8.85%'
    'G42' = '8'
    'D43' = '0.007701'
    'E43' = '4.10%'
    'G43' = '8'
    'E44' = '8.03%'
    'G44' = '8'
    'E45' = '-1.63%'
    'G45' = '8'
    'D46' = '0.00006998'
    'E46' = '6.21%'
    'G46' = '8'
    'D47' = '0.00000000750'
    'E47' = '-0.08%'
    'G47' = '8'
    'D48' = '0.05344'
    'E48' = '17.60%'
    'G48' = '8'
    'D49' = '0.004006'
    'G49' = '8'
    'D50' = '0.00002101'
    'E50' = '-0.08%'
    'G50' = '8'
    'D51' = '0.0002001'
    'E51' = '-0.08%'
    'G51' = '8'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
